$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.573
$ws.Range("C4").Value = 2.108
$ws.Range("D4").Value = 1.781
$ws.Range("E4").Value = 5.581
